# Effort tracking workbook update
# - Extend the "Fix: Bad specification..." note (row 25 / D25)
# - Bump the effort hours logged on 2012-10-15 (row 25 / B25): 1 -> 2.5
# - Append a new day of effort: 2012-10-16 (row 26) with its own note
# - Move the selection to reflect where the user ended up editing (A27)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing note text on row 25 (shared string reused in place) ---
$ws.Range("D25").Value = "Fix: Bad specification of ALL events - now timer events are still an OR condition. Implementation of waitForEventsTillTime by generalization of waitForEvent"

# --- Update effort hours already recorded for 2012-10-15 ---
$ws.Range("B25").Value = 2.5

# --- Add the new row 26 for 2012-10-16 ---
# Copy the date cell's format (style) from A25 so the new date cell keeps the
# same date number format without introducing a duplicate style entry.
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A26").Value = 41198
$ws.Range("B26").Value = 2.5
$ws.Range("D26").Value = "Code cleanup, suspendTillTime discarded. Implementation of enter/leaveCriticalSection"

# --- Restore view state: scroll a bit down and select the next empty row ---
$ws.Activate()
$ws.Range("A27").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
